$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 0.9070000052452087
$ws.Range("F2").Value = 0.8974641675854464
$ws.Range("J2").Value = 0.8139999999999999

$ws.Range("D3").Value = 0.9605000019073486
$ws.Range("F3").Value = 0.9588755856324831
$ws.Range("J3").Value = 0.921

$ws.Range("D4").Value = 0.9800000190734863
$ws.Range("F4").Value = 0.9795918367346939
$ws.Range("J4").Value = 0.96

$ws.Range("D5").Value = 0.9890000224113464
$ws.Range("F5").Value = 0.9888776541961577
$ws.Range("J5").Value = 0.978

$ws.Range("D6").Value = 0.996999979019165
$ws.Range("F6").Value = 0.9969939879759518
$ws.Range("H6").Value = 0.998995983935743
$ws.Range("J6").Value = 0.995
$ws.Range("L6").Value = 0.001

$ws.Range("D7").Value = 0.9909999966621399
$ws.Range("F7").Value = 0.9909182643794148
$ws.Range("J7").Value = 0.982

$ws.Range("D8").Value = 0.9934999942779541
$ws.Range("F8").Value = 0.9934574735782586
$ws.Range("J8").Value = 0.987

$ws.Range("D9").Value = 0.9990000128746033
$ws.Range("F9").Value = 0.9990009990009989
$ws.Range("H9").Value = 0.998003992015968
$ws.Range("L9").Value = 0.002

$ws.Range("D10").Value = 0.9990000128746033
$ws.Range("F10").Value = 0.998998998998999
$ws.Range("J10").Value = 0.998
